$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @{
    2 = @{ D=44315; L="Especial"; M=50; N=24000; O=24000; P=24000; Q="`$/caja 18 kilos"; R="Provincia de Melipilla"; S=1333; T=18 }
    3 = @{ D=44315; L="Primera"; M=50; N=20000; O=20000; P=20000; Q="`$/caja 18 kilos"; R="Provincia de Melipilla"; S=1111; T=18 }
    4 = @{ D=44291; L="Extra (doble especial)"; M=250; N=18000; O=18000; P=18000; Q="`$/caja 18 kilos"; R="Provincia de Melipilla"; S=1000; T=18 }
    5 = @{ D=44277; L="Especial"; M=200; N=15000; O=15000; P=15000; Q="`$/caja 18 kilos"; R="Provincia de Limarí"; S=833; T=18 }
    6 = @{ D=44292; L="Especial"; M=150; N=16000; O=16000; P=16000; Q="`$/caja 18 kilos"; R="Provincia de Melipilla"; S=889; T=18 }
    7 = @{ D=44292; L="Primera"; M=80; N=14000; O=14000; P=14000; Q="`$/caja 18 kilos"; R="Provincia de Melipilla"; S=778; T=18 }
    8 = @{ D=44279; L="Especial"; M=50; N=14000; O=14000; P=14000; Q="`$/caja 18 kilos"; R="Provincia de Melipilla"; S=778; T=18 }
    9 = @{ D=44279; L="Primera"; M=100; N=12000; O=12000; P=12000; Q="`$/caja 18 kilos"; R="Provincia de Melipilla"; S=667; T=18 }
    10 = @{ D=45022; L="Especial"; M=200; N=18000; O=18000; P=18000; Q="`$/caja 18 kilos"; R="Provincia de Melipilla"; S=1000; T=18 }
    11 = @{ D=44699; L="Especial"; M=150; N=22000; O=22000; P=22000; Q="`$/caja 18 kilos"; R="Provincia de Limarí"; S=1222; T=18 }
    12 = @{ D=44299; L="Especial"; M=170; N=18000; O=18000; P=18000; Q="`$/caja 18 kilos"; R="Provincia de Melipilla"; S=1000; T=18 }
    13 = @{ D=44299; L="Primera"; M=100; N=16000; O=16000; P=16000; Q="`$/caja 18 kilos"; R="Provincia de Melipilla"; S=889; T=18 }
    14 = @{ D=44630; L="Especial"; M=150; N=20000; O=20000; P=20000; Q="`$/caja 20 kilos"; R="Provincia de Limarí"; S=1000; T=20 }
    15 = @{ D=44698; L="Especial"; M=150; N=20000; O=20000; P=20000; Q="`$/caja 18 kilos"; R="Provincia de Limarí"; S=1111; T=18 }
    16 = @{ D=44698; L="Primera"; M=180; N=18000; O=18000; P=18000; Q="`$/caja 18 kilos"; R="Provincia de Limarí"; S=1000; T=18 }
    17 = @{ D=44985; L="Especial"; M=300; N=18000; O=18000; P=18000; Q="`$/caja 18 kilos"; R="Provincia de Limarí"; S=1000; T=18 }
    18 = @{ D=44985; L="Segunda"; M=150; N=12000; O=12000; P=12000; Q="`$/caja 18 kilos"; R="Provincia de Limarí"; S=667; T=18 }
    19 = @{ D=44350; L="Especial"; M=60; N=24000; O=24000; P=24000; Q="`$/caja 18 kilos"; R="Provincia de Limarí"; S=1333; T=18 }
    20 = @{ D=44300; L="Especial"; M=120; N=18000; O=18000; P=18000; Q="`$/caja 18 kilos"; R="Provincia de Melipilla"; S=1000; T=18 }
    21 = @{ D=44300; L="Primera"; M=100; N=16000; O=16000; P=16000; Q="`$/caja 18 kilos"; R="Provincia de Melipilla"; S=889; T=18 }
    22 = @{ D=44284; L="Especial"; M=120; N=13000; O=13000; P=13000; Q="`$/caja 18 kilos"; R="Provincia de Melipilla"; S=722; T=18 }
    23 = @{ D=44284; L="Extra (doble especial)"; M=100; N=15000; O=15000; P=15000; Q="`$/caja 18 kilos"; R="Provincia de Melipilla"; S=833; T=18 }
    24 = @{ D=44284; L="Primera"; M=50; N=12000; O=12000; P=12000; Q="`$/caja 18 kilos"; R="Provincia de Melipilla"; S=667; T=18 }
    25 = @{ D=44295; L="Segunda"; M=130; N=10000; O=10000; P=10000; Q="`$/caja 18 kilos"; R="Provincia de Melipilla"; S=556; T=18 }
    26 = @{ D=44301; L="Primera"; M=100; N=16000; O=16000; P=16000; Q="`$/caja 18 kilos"; R="Provincia de Melipilla"; S=889; T=18 }
    27 = @{ D=44224; L="Primera"; M=120; N=18000; O=18000; P=18000; Q="`$/caja 16 kilos"; R="Provincia de Limarí"; S=1125; T=16 }
    28 = @{ D=44645; L="Primera"; M=200; N=16000; O=16000; P=16000; Q="`$/caja 18 kilos"; R="Provincia de Limarí"; S=889; T=18 }
    29 = @{ D=44298; L="Extra (doble especial)"; M=160; N=20000; O=20000; P=20000; Q="`$/caja 18 kilos"; R="Provincia de Melipilla"; S=1111; T=18 }
    30 = @{ D=44271; L="Primera"; M=60; N=15000; O=15000; P=15000; Q="`$/caja 18 kilos"; R="Provincia de Melipilla"; S=833; T=18 }
    31 = @{ D=44252; L="Primera"; M=140; N=13000; O=13000; P=13000; Q="`$/caja 18 kilos"; R="Provincia de Melipilla"; S=722; T=18 }
    32 = @{ D=44274; L="Especial"; M=200; N=14000; O=14000; P=14000; Q="`$/caja 16 kilos"; R="Provincia de Melipilla"; S=875; T=16 }
    33 = @{ D=44274; L="Primera"; M=130; N=12000; O=12000; P=12000; Q="`$/caja 16 kilos"; R="Provincia de Melipilla"; S=750; T=16 }
    34 = @{ D=44258; L="Primera"; M=100; N=14000; O=14000; P=14000; Q="`$/caja 18 kilos"; R="Provincia de Limarí"; S=778; T=18 }
    35 = @{ D=44330; L="Primera"; M=50; N=23000; O=23000; P=23000; Q="`$/caja 18 kilos"; R="Provincia de Melipilla"; S=1278; T=18 }
    36 = @{ D=44267; L="Primera"; M=120; N=13000; O=13000; P=13000; Q="`$/caja 18 kilos"; R="Provincia de Melipilla"; S=722; T=18 }
    37 = @{ D=44222; L="Primera"; M=100; N=18000; O=18000; P=18000; Q="`$/caja 16 kilos"; R="Provincia de Limarí"; S=1125; T=16 }
    38 = @{ D=45083; L="Primera"; M=210; N=20000; O=20000; P=20000; Q="`$/caja 18 kilos"; R="Provincia de Melipilla"; S=1111; T=18 }
    39 = @{ D=45001; L="Especial"; M=150; N=13000; O=13000; P=13000; Q="`$/caja 18 kilos"; R="Provincia de Melipilla"; S=722; T=18 }
    40 = @{ D=45001; L="Primera"; M=100; N=11000; O=11000; P=11000; Q="`$/caja 18 kilos"; R="Provincia de Melipilla"; S=611; T=18 }
    41 = @{ D=44273; L="Especial"; M=40; N=15000; O=15000; P=15000; Q="`$/caja 16 kilos"; R="Provincia de Melipilla"; S=938; T=16 }
    42 = @{ D=44273; L="Primera"; M=50; N=13000; O=13000; P=13000; Q="`$/caja 16 kilos"; R="Provincia de Melipilla"; S=812; T=16 }
    43 = @{ D=44273; L="Segunda"; M=60; N=10000; O=10000; P=10000; Q="`$/caja 16 kilos"; R="Provincia de Melipilla"; S=625; T=16 }
    44 = @{ D=44309; L="Especial"; M=100; N=20000; O=20000; P=20000; Q="`$/caja 18 kilos"; R="Provincia de Melipilla"; S=1111; T=18 }
    45 = @{ D=44309; L="Primera"; M=60; N=18000; O=18000; P=18000; Q="`$/caja 18 kilos"; R="Provincia de Melipilla"; S=1000; T=18 }
}

foreach ($r in $rowData.Keys) {
    $d = $rowData[$r]
    $ws.Cells.Item($r, 4).Value = $d.D
    $ws.Cells.Item($r, 12).Value = $d.L
    $ws.Cells.Item($r, 13).Value = $d.M
    $ws.Cells.Item($r, 14).Value = $d.N
    $ws.Cells.Item($r, 15).Value = $d.O
    $ws.Cells.Item($r, 16).Value = $d.P
    $ws.Cells.Item($r, 17).Value = $d.Q
    $ws.Cells.Item($r, 18).Value = $d.R
    $ws.Cells.Item($r, 19).Value = $d.S
    $ws.Cells.Item($r, 20).Value = $d.T
}